# Update "想去人数" (want-to-go count, column F) values for a batch of
# conventions/events. The same set of events appears both on the "展览"
# sheet (sheet1) and on the "全部类型" sheet (sheet4, a superset of all
# event types), so each event's F-column value must be updated in both
# places, at its respective row.

$wb = $excel.ActiveWorkbook

# Locate the two worksheets by name (order-independent / robust).
$wsExhibit = $null
$wsAll = $null
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $s = $wb.Worksheets.Item($i)
    if ($s.Name -eq "展览") { $wsExhibit = $s }
    if ($s.Name -eq "全部类型") { $wsAll = $s }
}
if ($wsExhibit -eq $null) { $wsExhibit = $wb.Worksheets.Item(1) }
if ($wsAll -eq $null) { $wsAll = $wb.Worksheets.Item(4) }

# row on "展览" sheet -> row on "全部类型" sheet, plus old/new F values
# (F is column 6). Rows are matched 1:1 by event name between the sheets.
$updates = @(
    @{ R1 = 2;  R4 = 3;  Old = 97;    New = 101 },
    @{ R1 = 3;  R4 = 5;  Old = 12050; New = 12069 },
    @{ R1 = 4;  R4 = 6;  Old = 36;    New = 41 },
    @{ R1 = 5;  R4 = 7;  Old = 230;   New = 231 },
    @{ R1 = 6;  R4 = 9;  Old = 365;   New = 367 },
    @{ R1 = 8;  R4 = 11; Old = 11938; New = 11962 },
    @{ R1 = 9;  R4 = 12; Old = 502;   New = 503 },
    @{ R1 = 10; R4 = 13; Old = 1178;  New = 1180 },
    @{ R1 = 11; R4 = 14; Old = 110;   New = 111 },
    @{ R1 = 12; R4 = 15; Old = 583;   New = 588 },
    @{ R1 = 13; R4 = 16; Old = 1792;  New = 1796 },
    @{ R1 = 14; R4 = 18; Old = 5910;  New = 5917 },
    @{ R1 = 16; R4 = 20; Old = 3556;  New = 3557 },
    @{ R1 = 17; R4 = 21; Old = 198;   New = 200 },
    @{ R1 = 18; R4 = 22; Old = 30;    New = 32 }
)

$mismatches = 0
foreach ($u in $updates) {
    $c1 = $wsExhibit.Cells.Item($u.R1, 6)
    $c4 = $wsAll.Cells.Item($u.R4, 6)

    if ($c1.Value() -ne $u.Old) { $mismatches++ }
    if ($c4.Value() -ne $u.Old) { $mismatches++ }

    $c1.Value = $u.New
    $c4.Value = $u.New
}

Write-Host "Updated" $updates.Count "events on both sheets. Mismatches vs expected old values:" $mismatches
